$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H76").Value = 3165.8
$ws.Range("I76").Value = 3180
$ws.Range("J76").Value = 3073.5
$ws.Range("K76").Value = 3180
$ws.Range("L76").Value = 3073.5
$ws.Range("M76").Value = -2865
$ws.Range("N76").Value = -3703.5
$ws.Range("H79").Value = 3165.8
$ws.Range("I79").Value = 3180
$ws.Range("J79").Value = 3073.5
$ws.Range("K79").Value = 3180
$ws.Range("L79").Value = 3073.5
$ws.Range("M79").Value = -2088
$ws.Range("N79").Value = -5257.5
$ws.Range("H118").Value = 2070.9333
$ws.Range("I118").Value = 1506.9
$ws.Range("J118").Value = 3199
$ws.Range("K118").Value = 4520.700000000001
$ws.Range("L118").Value = 9597
$ws.Range("M118").Value = -2863.700000000001
$ws.Range("N118").Value = -12911
$ws.Range("H127").Value = 1084.8948
$ws.Range("I127").Value = 414.44446
$ws.Range("J127").Value = 1688.3
$ws.Range("K127").Value = 1243.33338
$ws.Range("L127").Value = 5064.9
$ws.Range("M127").Value = 3716.66662
$ws.Range("N127").Value = -14984.9
$ws.Range("H129").Value = 1351.4036
$ws.Range("I129").Value = 446.55
$ws.Range("J129").Value = 1840.5135
$ws.Range("K129").Value = 1339.65
$ws.Range("L129").Value = 5521.5405
$ws.Range("M129").Value = 3660.35
$ws.Range("N129").Value = -15521.5405
$ws.Range("H132").Value = 2373.2964
$ws.Range("I132").Value = 2401.7112
$ws.Range("J132").Value = 2231.2222
$ws.Range("K132").Value = 7205.133600000001
$ws.Range("L132").Value = 6693.6666
$ws.Range("M132").Value = -4675.133600000001
$ws.Range("N132").Value = -11753.6666
$ws.Range("H133").Value = 79780
$ws.Range("J133").Value = 79780
$ws.Range("L133").Value = 79780
$ws.Range("N133").Value = -89900
$ws.Range("H137").Value = 1204.1754
$ws.Range("I137").Value = 1031.9756
$ws.Range("J137").Value = 1645.4375
$ws.Range("K137").Value = 3095.9268
$ws.Range("L137").Value = 4936.3125
$ws.Range("M137").Value = -545.9268000000002
$ws.Range("N137").Value = -10036.3125
$ws.Range("H138").Value = 2049.6155
$ws.Range("I138").Value = 1361.4314
$ws.Range("J138").Value = 2927.05
$ws.Range("K138").Value = 4084.2942
$ws.Range("L138").Value = 8781.150000000001
$ws.Range("M138").Value = 1055.7058
$ws.Range("N138").Value = -19061.15

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 1103528.4
$ws.Range("I32").Value = 1280093.8
$ws.Range("K32").Value = 1280093.8
$ws.Range("M32").Value = -1279806.8
$ws.Range("H45").Value = 1970.0476
$ws.Range("I45").Value = 1931.8
$ws.Range("K45").Value = 1931.8
$ws.Range("M45").Value = -1554.8
$ws.Range("H61").Value = 2127.4856
$ws.Range("I61").Value = 1851.5957
$ws.Range("J61").Value = 2691.261
$ws.Range("K61").Value = 1851.5957
$ws.Range("L61").Value = 2691.261
$ws.Range("M61").Value = -1639.5957
$ws.Range("N61").Value = -3115.261
$ws.Range("H74").Value = 987.6731
$ws.Range("I74").Value = 740.125
$ws.Range("J74").Value = 1383.75
$ws.Range("K74").Value = 740.125
$ws.Range("L74").Value = 1383.75
$ws.Range("M74").Value = 133.875
$ws.Range("N74").Value = -3131.75
$ws.Range("H77").Value = 987.6731
$ws.Range("I77").Value = 740.125
$ws.Range("J77").Value = 1383.75
$ws.Range("K77").Value = 3700.625
$ws.Range("L77").Value = 6918.75
$ws.Range("M77").Value = 667.375
$ws.Range("N77").Value = -15654.75
$ws.Range("H122").Value = 144107.14
$ws.Range("I122").Value = 200747.2
$ws.Range("J122").Value = 2507
$ws.Range("K122").Value = 602241.6000000001
$ws.Range("L122").Value = 7521
$ws.Range("M122").Value = -599791.6000000001
$ws.Range("N122").Value = -12421
$ws.Range("H132").Value = 3686.4783
$ws.Range("I132").Value = 2976
$ws.Range("J132").Value = 5490
$ws.Range("K132").Value = 8928
$ws.Range("L132").Value = 16470
$ws.Range("M132").Value = -6398
$ws.Range("N132").Value = -21530
$ws.Range("H136").Value = 2127.4856
$ws.Range("I136").Value = 1851.5957
$ws.Range("J136").Value = 2691.261
$ws.Range("K136").Value = 5554.7871
$ws.Range("L136").Value = 8073.782999999999
$ws.Range("M136").Value = -3004.7871
$ws.Range("N136").Value = -13173.783

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H107").Value = 144404.86
$ws.Range("I107").Value = 250958.5
$ws.Range("J107").Value = 2333.3333
$ws.Range("K107").Value = 250958.5
$ws.Range("L107").Value = 2333.3333
$ws.Range("M107").Value = -249038.5
$ws.Range("N107").Value = -6173.3333

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 5893.2646
$ws.Range("I31").Value = 2003.7273
$ws.Range("J31").Value = 7753.478
$ws.Range("K31").Value = 2003.7273
$ws.Range("L31").Value = 7753.478
$ws.Range("M31").Value = -1708.7273
$ws.Range("N31").Value = -8343.477999999999
$ws.Range("H34").Value = 5893.2646
$ws.Range("I34").Value = 2003.7273
$ws.Range("J34").Value = 7753.478
$ws.Range("K34").Value = 2003.7273
$ws.Range("L34").Value = 7753.478
$ws.Range("M34").Value = -1801.7273
$ws.Range("N34").Value = -8157.478
$ws.Range("H58").Value = 1164.5853
$ws.Range("I58").Value = 887.3913
$ws.Range("K58").Value = 887.3913
$ws.Range("M58").Value = -684.3913
$ws.Range("H132").Value = 4275254
$ws.Range("I132").Value = 1922.7273
$ws.Range("J132").Value = 9805448
$ws.Range("K132").Value = 5768.1819
$ws.Range("L132").Value = 29416344
$ws.Range("M132").Value = -3238.1819
$ws.Range("N132").Value = -29421404
$ws.Range("H136").Value = 1164.5853
$ws.Range("I136").Value = 887.3913
$ws.Range("K136").Value = 2662.1739
$ws.Range("M136").Value = -112.1738999999998

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H38").Value = 123.14286
$ws.Range("I38").Value = 20
$ws.Range("J38").Value = 260.66666
$ws.Range("K38").Value = 60
$ws.Range("L38").Value = 781.9999799999999
$ws.Range("M38").Value = 287
$ws.Range("N38").Value = -1475.99998
$ws.Range("H131").Value = 2743.4492
$ws.Range("J131").Value = 2917.0312
$ws.Range("L131").Value = 8751.0936
$ws.Range("N131").Value = -18831.0936

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H97").Value = 583.9048
$ws.Range("I97").Value = 603.125
$ws.Range("J97").Value = 522.4
$ws.Range("K97").Value = 603.125
$ws.Range("L97").Value = 522.4
$ws.Range("M97").Value = -107.125
$ws.Range("N97").Value = -1514.4
$ws.Range("H122").Value = 3400
$ws.Range("I122").Value = 2800
$ws.Range("J122").Value = 3600
$ws.Range("K122").Value = 8400
$ws.Range("L122").Value = 10800
$ws.Range("M122").Value = -5950
$ws.Range("N122").Value = -15700
$ws.Range("H125").Value = 40326
$ws.Range("J125").Value = 40326
$ws.Range("L125").Value = 40326
$ws.Range("N125").Value = -45246

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H16").Value = 15876456
$ws.Range("I16").Value = 3871.125
$ws.Range("K16").Value = 3871.125
$ws.Range("M16").Value = -3701.125
$ws.Range("H62").Value = 100000
$ws.Range("J62").Value = 100000
$ws.Range("L62").Value = 100000
$ws.Range("N62").Value = -101248
$ws.Range("H65").Value = 100000
$ws.Range("J65").Value = 100000
$ws.Range("L65").Value = 300000
$ws.Range("N65").Value = -306240
$ws.Range("H122").Value = 3954.6155
$ws.Range("I122").Value = 2663.3333
$ws.Range("J122").Value = 4342
$ws.Range("K122").Value = 7989.999899999999
$ws.Range("L122").Value = 13026
$ws.Range("M122").Value = -5539.999899999999
$ws.Range("N122").Value = -17926
$ws.Range("H132").Value = 2081.0679
$ws.Range("I132").Value = 2061.3235
$ws.Range("J132").Value = 2107.92
$ws.Range("K132").Value = 6183.970499999999
$ws.Range("L132").Value = 6323.76
$ws.Range("M132").Value = -3653.970499999999
$ws.Range("N132").Value = -11383.76
$ws.Range("H136").Value = 7248281.5
$ws.Range("I136").Value = 1951.5294
$ws.Range("K136").Value = 5854.5882
$ws.Range("M136").Value = -3304.5882

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H122").Value = 3079.7
$ws.Range("I122").Value = 3079.7
$ws.Range("K122").Value = 9239.099999999999
$ws.Range("M122").Value = -6789.099999999999
$ws.Range("H132").Value = 2193940.2
$ws.Range("I132").Value = 926.8421
$ws.Range("J132").Value = 8772981
$ws.Range("K132").Value = 2780.5263
$ws.Range("L132").Value = 26318943
$ws.Range("M132").Value = -250.5263
$ws.Range("N132").Value = -26324003
$ws.Range("H136").Value = 2804.2407
$ws.Range("I136").Value = 2532.425
$ws.Range("J136").Value = 3580.8572
$ws.Range("K136").Value = 7597.275000000001
$ws.Range("L136").Value = 10742.5716
$ws.Range("M136").Value = -5047.275000000001
$ws.Range("N136").Value = -15842.5716

Write-Host "All changes applied."